$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.395.39'
$ws.Range("E2").Value = '  -0.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.507.15'
$ws.Range("E3").Value = '  -1.43%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.25'
$ws.Range("E5").Value = '  +0.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.00'
$ws.Range("E6").Value = '  -3.57%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.560'
$ws.Range("E7").Value = '  -1.83%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  -3.49%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.93'
$ws.Range("E10").Value = '  -3.06%  '

# Row 11
$ws.Range("E11").Value = '  -1.24%  '

# Row 12
$ws.Range("E12").Value = '  +0.40%  '

# Row 13
$ws.Range("E13").Value = '  -3.58%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.891.53'
$ws.Range("E14").Value = '  -1.59%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.527.81'
$ws.Range("E15").Value = '  -0.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.13'
$ws.Range("E16").Value = '  -5.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.374.20'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.53'
$ws.Range("E19").Value = '  -4.23%  '

# Row 20
$ws.Range("E20").Value = '  -2.45%  '

# Row 21
$ws.Range("E21").Value = '  -3.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.29'
$ws.Range("E22").Value = '  -1.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.35'
$ws.Range("E23").Value = '  -2.05%  '

# Row 24
$ws.Range("E24").Value = '  -3.10%  '

# Row 26
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.15'
$ws.Range("E27").Value = '  -4.92%  '

# Row 28
$ws.Range("E28").Value = '  -4.76%  '

# Row 29
$ws.Range("E29").Value = '  -1.64%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.40'
$ws.Range("E30").Value = '  -7.97%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.19'
$ws.Range("E31").Value = '  -1.17%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.77'
$ws.Range("E32").Value = '  +1.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  -0.53%  '

# Row 34
$ws.Range("E34").Value = '  +0.66%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0776'
$ws.Range("E35").Value = '  -3.45%  '

# Row 36
$ws.Range("E36").Value = '  -3.91%  '

# Row 37
$ws.Range("E37").Value = '  -6.40%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.33'
$ws.Range("E38").Value = '  -5.80%  '

# Row 39
$ws.Range("E39").Value = '  -5.00%  '

# Row 40
$ws.Range("E40").Value = '  -1.25%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.11'
$ws.Range("E41").Value = '  -1.89%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.29'
$ws.Range("E42").Value = '  -4.59%  '

# Row 43
$ws.Range("E43").Value = '  -0.12%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0292'
$ws.Range("E44").Value = '  -1.98%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.995.36'
$ws.Range("E45").Value = '  +1.07%  '

# Row 46
$ws.Range("E46").Value = '  -4.08%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.84'
$ws.Range("E47").Value = '  -0.61%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.745.30'
$ws.Range("E48").Value = '  -1.66%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '78.18'
$ws.Range("E49").Value = '  -3.67%  '

# Row 50
$ws.Range("E50").Value = '  -3.64%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.79'
$ws.Range("E51").Value = '  -3.29%  '
